# The underlying reference data on the "ID_HeatingTechnology" sheet was
# reworked (infrastructure availability data developed). The cell values
# themselves did not change, but re-saving the refreshed workbook causes
# Excel to recompute the best-fit column widths for the two data columns
# and leaves the cursor on the last-touched cell (B19) instead of the
# previous selection (D18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-fit both data columns ("id_heating_technology" / "name") to their
# (refreshed) content, matching the new best-fit widths.
$ws.Columns.Item(1).ColumnWidth = 23.1665
$ws.Columns.Item(2).ColumnWidth = 50.5

# Leave the selection where the edit was made.
$ws.Range("B19").Select() | Out-Null
